# Commit: "modified born position of city"
# Update the RelivePos (born position) value for the first scene row (villageScene/City)
# from "0,0,0" to "20,0,-137".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "20,0,-137"
